$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to the new board name "Board_1452" (adds a new shared string entry)
$ws.Range("A2").Value = "Board_1452"

# Update the selection to match the diff (A2 active, sqref A2)
$ws.Range("A2").Select()
